$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.411.78"
$ws.Range("E2").Value = "  +0.42%  "

$ws.Range("D3").Value = "3.745.98"
$ws.Range("E3").Value = "  -0.96%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'592.63"
$ws.Range("E5").Value = "  -0.73%  "

$ws.Range("D6").Value = "'165.86"
$ws.Range("E6").Value = "  -1.64%  "

$ws.Range("D7").Value = "3.742.69"
$ws.Range("E7").Value = "  -1.06%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("E9").Value = "  -1.41%  "

$ws.Range("E10").Value = "  -2.90%  "

$ws.Range("D11").Value = "'6.41"
$ws.Range("E11").Value = "  -1.57%  "

$ws.Range("D12").Value = "'0.448"
$ws.Range("E12").Value = "  -1.00%  "

$ws.Range("E13").Value = "  -6.82%  "

$ws.Range("D14").Value = "'35.92"
$ws.Range("E14").Value = "  -2.32%  "

$ws.Range("D15").Value = "4.372.95"
$ws.Range("E15").Value = "  -1.04%  "

$ws.Range("D16").Value = "3.741.71"
$ws.Range("E16").Value = "  -1.09%  "

$ws.Range("D17").Value = "68.360.28"
$ws.Range("E17").Value = "  +0.43%  "

$ws.Range("D18").Value = "'17.90"
$ws.Range("E18").Value = "  -4.08%  "

$ws.Range("E19").Value = "  +0.03%  "

$ws.Range("D20").Value = "'6.95"
$ws.Range("E20").Value = "  -2.79%  "

$ws.Range("E21").Value = "  +1.07%  "

$ws.Range("D22").Value = "'463.21"
$ws.Range("E22").Value = "  -1.01%  "

$ws.Range("D23").Value = "'0.694"
$ws.Range("E23").Value = "  -3.35%  "

$ws.Range("D24").Value = "'83.87"
$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("D25").Value = "'0.0000145"
$ws.Range("E25").Value = "  -1.17%  "

$ws.Range("D26").Value = "'2.17"
$ws.Range("E26").Value = "  -3.25%  "

$ws.Range("D27").Value = "'11.88"
$ws.Range("E27").Value = "  -2.45%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'10.03"
$ws.Range("E28").Value = "  -3.73%  "

$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.12%  "

$ws.Range("D30").Value = "3.890.69"
$ws.Range("E30").Value = "  -1.04%  "

$ws.Range("D31").Value = "'2.76"
$ws.Range("E31").Value = "  -5.55%  "

$ws.Range("D32").Value = "'7.28"
$ws.Range("E32").Value = "  -3.89%  "

$ws.Range("D33").Value = "'29.84"
$ws.Range("E33").Value = "  -2.02%  "

$ws.Range("E34").Value = "  -3.67%  "

$ws.Range("D35").Value = "'1.00"

$ws.Range("D36").Value = "'9.11"
$ws.Range("E36").Value = "  -1.05%  "

$ws.Range("D37").Value = "3.701.30"
$ws.Range("E37").Value = "  -1.07%  "

$ws.Range("D38").Value = "'0.100"
$ws.Range("E38").Value = "  -3.65%  "

$ws.Range("D39").Value = "'3.38"
$ws.Range("E39").Value = "  -9.07%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.137"
$ws.Range("E40").Value = "  -1.31%  "

$ws.Range("B41").Value = "Mantle"
$ws.Range("C41").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D41").Value = "'0.995"
$ws.Range("E41").Value = "  -1.20%  "

$ws.Range("D42").Value = "'5.76"
$ws.Range("E42").Value = "  -1.60%  "

$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  -0.10%  "

$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("D45").Value = "'0.300"
$ws.Range("E45").Value = "  -4.11%  "

$ws.Range("D46").Value = "'43.40"
$ws.Range("E46").Value = "  +8.61%  "

$ws.Range("D47").Value = "'46.55"
$ws.Range("E47").Value = "  +2.13%  "

$ws.Range("D48").Value = "'1.91"
$ws.Range("E48").Value = "  -1.74%  "

$ws.Range("D49").Value = "'8.45"
$ws.Range("E49").Value = "  -2.56%  "

$ws.Range("D50").Value = "'144.52"
$ws.Range("E50").Value = "  +0.32%  "

$ws.Range("D51").Value = "'387.80"
$ws.Range("E51").Value = "  -4.38%  "
